$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.731.93"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.861.28"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'1.021"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "'320.69"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.4374"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.3803"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.07439"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'0.8824"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'21.61"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "1.860.99"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'6.745"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "'5.492"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'0.07136"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "'86.55"
$ws.Range("E16").Value = "  +4.80%  "
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'0.000009070"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "'15.48"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "27.730.04"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'5.291"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").Value = "2.094.65"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'2.041"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").Value = "'157.53"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "'2.003"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").Value = "'5.359"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "'121.32"
$ws.Range("E30").Value = "  +3.99%  "
$ws.Range("D31").Value = "'0.09052"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "'1.214"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").Value = "'0.7661"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'3.030"
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("D35").Value = "'4.565"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").Value = "'1.020"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "'0.01982"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "'0.05294"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'2.871"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").Value = "'0.5194"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'6.953"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "'8.703"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").Value = "'10.80"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").Value = "'110.12"
$ws.Range("D47").Value = "'1.716"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").Value = "'1.022"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'0.06501"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'0.4712"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "'1.883"
$ws.Range("E51").Value = "  +1.57%  "
